# edit.ps1 - Word COM-interop script applying the commit's changes:
#   1. Trim the "1.2.4 Post-appointment (if applicable): ..." run down to
#      "1.2.4 Post-appointment (if applicable)" and leave a "_GoBack"
#      bookmark (empty span) at the end of that paragraph - this is what
#      Word stamps at the location of the most recent edit.
#   2. Mark the run containing "4.1.1 Create account" with a
#      <w:lastRenderedPageBreak/> (Word inserts these while paginating;
#      here it lands right before that heading).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Shorten the 1.2.4 paragraph text and drop a _GoBack bookmark after it
# ---------------------------------------------------------------------
$oldText = "1.2.4 Post-appointment (if applicable): Receive confirmation sms or email, prompt customer to share experience on review sites"
$newText = "1.2.4 Post-appointment (if applicable)"

$findRange = $d.Content
$found = $findRange.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

if ($found) {
    foreach ($para in $d.Paragraphs) {
        if ($para.Range.Text -eq ($newText + [char]13)) {
            # Position right after the last character of the paragraph's
            # text, i.e. right before its trailing paragraph mark. Adding
            # a bookmark straight at that (empty) position is flaky, so
            # nudge it into place: insert a throwaway character there,
            # anchor the bookmark to it, then delete the character again -
            # the (now empty) bookmark stays put.
            $endPos = $para.Range.End - 1
            $anchor = $d.Range($endPos, $endPos)
            $anchor.InsertAfter("x")
            $d.Bookmarks.Add("_GoBack", $anchor) | Out-Null
            $anchor.Text = ""
            break
        }
    }
}

# ---------------------------------------------------------------------
# 2) Add <w:lastRenderedPageBreak/> right before "4.1.1 Create account"
# ---------------------------------------------------------------------
$target = "4.1.1 Create account"
foreach ($para in @($d.Paragraphs)) {
    if ($para.Range.Text -eq ($target + [char]13)) {
        $runRange = $d.Range($para.Range.Start, $para.Range.End - 1)
        $packageXml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p><w:r w:rsidRPr="18804DFF"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t>$target</w:t></w:r></w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@
        $runRange.InsertXML($packageXml)
        break
    }
}

Write-Host "Done"
